$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - (Intercept)
$ws.Range("B2").Value = 0.800723
$ws.Range("D2").Value = 2.338088
$ws.Range("E2").Value = 0.132947

# Row 3 - household_group_collapsed
$ws.Range("B3").Value = 10.704165
$ws.Range("D3").Value = 15.627918
$ws.Range("E3").Value = 0.000006

# Row 4 - Residuals
$ws.Range("B4").Value = 16.096058
$ws.Range("C4").Value = 47

# Row 5 - SM-Control
$ws.Range("G5").Value = 0.067159
$ws.Range("H5").Value = -0.719893
$ws.Range("I5").Value = 0.85421
$ws.Range("J5").Value = 0.976776

# Row 6 - SM + Traps-Control
$ws.Range("G6").Value = 0.991187
$ws.Range("H6").Value = 0.235788
$ws.Range("I6").Value = 1.746585
$ws.Range("J6").Value = 0.007312

# Row 7 - SM + Traps-SM
$ws.Range("G7").Value = 0.924028
$ws.Range("H7").Value = 0.491411
$ws.Range("I7").Value = 1.356645
$ws.Range("J7").Value = 0.000014
